# Commit: Wed, Jun 17, 2020 1:05:42 PM
#
# The underlying edit is a table-style change on the table that lives on
# slide 16 (the "PLENARY- COMPLETE THE MISSING GAPS" slide): its table
# style id goes from {06FCCF8C-1D8A-4F04-B829-4A845869C854} (the style
# embedded in ppt/tableStyles.xml) to {D12CDB09-35F2-45B8-958B-55FB48062636}
# (PowerPoint's built-in default table style, "Medium Style 2 - Accent 1").

$p = $ppt.ActivePresentation
$targetStyleId = "{D12CDB09-35F2-45B8-958B-55FB48062636}"

$slide = $p.Slides.Item(16)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $table = $shape.Table
        # Table styles are reassigned through ApplyStyle, not by setting
        # the Style property directly.
        $table.ApplyStyle($targetStyleId)
    }
}
